# Weekly fruit/vegetable price update:
# Insert a new week's record as row 631 (pushing all subsequent rows down by
# one), matching the "Fruta / hortaliza, semanal" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 631:715 down to 632:716 by inserting a new row at 631.
$ws.Rows.Item(631).Insert()

# Populate the newly inserted row 631 with this week's record.
$ws.Cells.Item(631, 1).Value  = 4
$ws.Cells.Item(631, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(631, 3).Value  = "Los Lagos"
$ws.Cells.Item(631, 4).Value  = 44776
$ws.Cells.Item(631, 5).Value  = 10
$ws.Cells.Item(631, 6).Value  = 100112020
$ws.Cells.Item(631, 7).Value  = "Tomate"
$ws.Cells.Item(631, 8).Value  = "Larga vida"
$ws.Cells.Item(631, 9).Value  = "Primera"
$ws.Cells.Item(631, 10).Value = 120
$ws.Cells.Item(631, 11).Value = 16000
$ws.Cells.Item(631, 12).Value = 16000
$ws.Cells.Item(631, 13).Value = 16000
$ws.Cells.Item(631, 14).Value = "$/bandeja 20 kilos"
$ws.Cells.Item(631, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(631, 16).Value = 800
$ws.Cells.Item(631, 17).Value = 20
$ws.Cells.Item(631, 18).Value = "Hortaliza"
